$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $formulaExpr) {
    # Write as a formula (Excel formula syntax) producing the literal text,
    # then convert the cell to a static value via copy/paste-values. This
    # guarantees the cell ends up as plain text (matching the source
    # t="inlineStr" cells) instead of Excel auto-coercing numeric-looking
    # text ("5.40") into a Number and silently dropping trailing zeros, and
    # without leaving any new/unused number-format style behind.
    $ws.Range($cellRef).Formula = "=" + $formulaExpr
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
}

Set-TextValue 'D2' '"37.053.82"'
$ws.Range('E2').Value = '  -0.26%  '
Set-TextValue 'D3' '"2.048.05"'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  +0.11%  '
Set-TextValue 'D5' '"246.14"'
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('E6').Value = '  -2.06%  '
Set-TextValue 'D7' '"58.68"'
$ws.Range('E7').Value = '  -4.90%  '
$ws.Range('E9').Value = '  -2.18%  '
$ws.Range('E10').Value = '  -2.48%  '
Set-TextValue 'D11' '"0.111"'
$ws.Range('E11').Value = '  +2.26%  '
Set-TextValue 'D12' '"15.39"'
$ws.Range('E12').Value = '  -5.87%  '
Set-TextValue 'D13' '"0.892"'
$ws.Range('E13').Value = '  +7.94%  '
Set-TextValue 'D14' '"2.341.70"'
$ws.Range('E14').Value = '  -0.61%  '
Set-TextValue 'D15' '"5.73"'
$ws.Range('E15').Value = '  +0.03%  '
Set-TextValue 'D16' '"2.012.49"'
$ws.Range('E16').Value = '  -2.16%  '
Set-TextValue 'D17' '"18.35"'
$ws.Range('E17').Value = '  +1.82%  '
Set-TextValue 'D18' '"37.042.29"'
$ws.Range('E18').Value = '  -0.29%  '
Set-TextValue 'D19' '"73.82"'
$ws.Range('E19').Value = '  -1.94%  '
Set-TextValue 'D20' '"0.0" & UNICHAR(8323) & "0886"'
$ws.Range('E20').Value = '  -2.10%  '
Set-TextValue 'D21' '"5.40"'
$ws.Range('E21').Value = '  -0.70%  '
Set-TextValue 'D22' '"240.20"'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('E23').Value = '  -0.02%  '
Set-TextValue 'D24' '"2.44"'
$ws.Range('E24').Value = '  +1.17%  '
Set-TextValue 'D25' '"9.68"'
$ws.Range('E25').Value = '  +2.47%  '
Set-TextValue 'D26' '"168.54"'
$ws.Range('E26').Value = '  -0.57%  '
Set-TextValue 'D27' '"2.15"'
$ws.Range('E27').Value = '  -2.96%  '
$ws.Range('E28').Value = '  -0.32%  '
Set-TextValue 'D29' '"5.56"'
$ws.Range('E29').Value = '  +15.20%  '
$ws.Range('E30').Value = '  -1.05%  '
$ws.Range('E31').Value = '  -2.37%  '
$ws.Range('E32').Value = '  +4.17%  '
Set-TextValue 'D33' '"0.0614"'
$ws.Range('E33').Value = '  -1.32%  '
$ws.Range('E34').Value = '  +0.15%  '
Set-TextValue 'D35' '"1.85"'
$ws.Range('E35').Value = '  +6.27%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D36' '"2.26"'
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D37' '"0.0850"'
$ws.Range('E37').Value = '  -5.33%  '
$ws.Range('E38').Value = '  -3.59%  '
Set-TextValue 'D39' '"5.26"'
$ws.Range('E39').Value = '  -1.24%  '
$ws.Range('E40').Value = '  -2.32%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D41' '"0.0223"'
$ws.Range('E41').Value = '  -0.75%  '
$ws.Range('B42').Value = 'Cronos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D42' '"0.0972"'
$ws.Range('E42').Value = '  -10.79%  '
Set-TextValue 'D43' '"1.15"'
Set-TextValue 'D44' '"97.59"'
$ws.Range('E44').Value = '  -0.58%  '
Set-TextValue 'D45' '"17.03"'
$ws.Range('E45').Value = '  -6.43%  '
Set-TextValue 'D46' '"2.39"'
$ws.Range('E46').Value = '  -4.05%  '
Set-TextValue 'D47' '"1.299.53"'
$ws.Range('E47').Value = '  +0.16%  '
Set-TextValue 'D48' '"2.87"'
$ws.Range('E48').Value = '  -0.23%  '
Set-TextValue 'D49' '"6.75"'
$ws.Range('E49').Value = '  -1.65%  '
$ws.Range('E50').Value = '  +2.70%  '
Set-TextValue 'D51' '"2.227.36"'
$ws.Range('E51').Value = '  -0.67%  '

$excel.CutCopyMode = $false
